# Update "想去人数" (number of people interested) figures to the latest
# scraped values, as produced by the site generator run at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 45
$ws1.Range("F5").Value = 2551
$ws1.Range("F6").Value = 235
$ws1.Range("F7").Value = 380

# --- Sheet "全部类型" (All types, combined view) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 45
$ws4.Range("F5").Value = 2551
$ws4.Range("F6").Value = 235
$ws4.Range("F9").Value = 380
